$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Randomly transferred student quotas between courses
$ws.Range("C2").Value = 80
$ws.Range("C3").Value = 60
$ws.Range("C6").Value = 120
$ws.Range("C10").Value = 28

# Update the active selection to reflect the last-touched cell
$ws.Range("C11").Select()
